# Weekly fruit/vegetable price data refresh: a new week's data point is
# inserted at row 293, pushing all existing data rows down by two rows
# (the engine's Insert default duplicates the preceding row's formatting,
# which matches the original workbook's date-style column D).
#
# After the insert:
#   - row 293 becomes a brand-new record (date 2022-06-10 / serial 44722)
#   - row 294 becomes a brand-new record (also date serial 44722)
#   - the old rows 293..365 now live at rows 295..367

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("293:294").Insert()

# --- New row 293 ---
$ws.Cells.Item(293, 1).Value = 10
$ws.Cells.Item(293, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(293, 3).Value = "La Araucanía"
$ws.Cells.Item(293, 4).Value = 44722
$ws.Cells.Item(293, 5).Value = 9
$ws.Cells.Item(293, 6).Value = 100112037
$ws.Cells.Item(293, 7).Value = "Cebollín"
$ws.Cells.Item(293, 8).Value = "Sin especificar"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 30
$ws.Cells.Item(293, 11).Value = 10000
$ws.Cells.Item(293, 12).Value = 10000
$ws.Cells.Item(293, 13).Value = 10000
$ws.Cells.Item(293, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(293, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(293, 16).Value = 833
$ws.Cells.Item(293, 17).Value = 12
$ws.Cells.Item(293, 18).Value = "Hortaliza"

# --- New row 294 ---
$ws.Cells.Item(294, 1).Value = 10
$ws.Cells.Item(294, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(294, 3).Value = "La Araucanía"
$ws.Cells.Item(294, 4).Value = 44722
$ws.Cells.Item(294, 5).Value = 9
$ws.Cells.Item(294, 6).Value = 100112037
$ws.Cells.Item(294, 7).Value = "Cebollín"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 30
$ws.Cells.Item(294, 11).Value = 7000
$ws.Cells.Item(294, 12).Value = 7000
$ws.Cells.Item(294, 13).Value = 7000
$ws.Cells.Item(294, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(294, 15).Value = "Región Metropolitana"
$ws.Cells.Item(294, 16).Value = 583
$ws.Cells.Item(294, 17).Value = 12
$ws.Cells.Item(294, 18).Value = "Hortaliza"
